$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.112.94'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.318.19'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '253.75'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').Value = '0.643'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '76.18'
$ws.Range('E7').Value = '  +5.94%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('D10').Value = '39.77'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('D12').Value = '7.64'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').Value = '2.659.78'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '15.53'
$ws.Range('E15').Value = '  +4.23%  '
$ws.Range('D16').Value = '0.886'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '2.312.05'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').Value = '43.084.46'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').Value = '6.32'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '73.11'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = '239.08'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '2.25'
$ws.Range('E23').Value = '  +6.00%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '11.57'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = '2.44'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.15'
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '21.33'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '167.76'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '6.39'
$ws.Range('E31').Value = '  +2.78%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0844'
$ws.Range('E32').Value = '  +8.59%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '0.129'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '30.67'
$ws.Range('E34').Value = '  +4.65%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '0.129'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  +10.24%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = '4.86'
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0315'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '13.86'
$ws.Range('E39').Value = '  +11.69%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '2.36'
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '5.92'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.219'
$ws.Range('E42').Value = '  +8.48%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '9.24'
$ws.Range('E43').Value = '  +3.12%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = '62.79'
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '4.93'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '105.40'
$ws.Range('E46').Value = '  +10.62%  '
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '1.19'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('B49').Value = 'BinanceUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.19'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '4.38'
$ws.Range('E51').Value = '  -1.07%  '
